$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'58.315.87"
$ws.Range("E2").Value = "'  +0.67%  "

# Row 3
$ws.Range("D3").Value = "'2.523.86"
$ws.Range("E3").Value = "'  +2.48%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'521.02"
$ws.Range("E5").Value = "'  +0.58%  "

# Row 6
$ws.Range("D6").Value = "'132.07"
$ws.Range("E6").Value = "'  +1.09%  "

# Row 7
$ws.Range("E7").Value = "'  -0.12%  "

# Row 8
$ws.Range("D8").Value = "'0.556"
$ws.Range("E8").Value = "'  -0.24%  "

# Row 9
$ws.Range("D9").Value = "'2.521.31"
$ws.Range("E9").Value = "'  +2.44%  "

# Row 10
$ws.Range("D10").Value = "'0.0974"
$ws.Range("E10").Value = "'  -1.57%  "

# Row 11
$ws.Range("D11").Value = "'0.157"
$ws.Range("E11").Value = "'  -0.07%  "

# Row 12
$ws.Range("D12").Value = "'5.22"
$ws.Range("E12").Value = "'  -2.16%  "

# Row 13
$ws.Range("E13").Value = "'  -2.11%  "

# Row 14
$ws.Range("D14").Value = "'2.963.94"
$ws.Range("E14").Value = "'  +2.18%  "

# Row 15
$ws.Range("D15").Value = "'58.283.46"
$ws.Range("E15").Value = "'  +0.74%  "

# Row 16
$ws.Range("D16").Value = "'22.25"
$ws.Range("E16").Value = "'  -0.07%  "

# Row 17
$ws.Range("E17").Value = "'  -0.41%  "

# Row 18
$ws.Range("D18").Value = "'2.515.58"
$ws.Range("E18").Value = "'  +2.14%  "

# Row 19
$ws.Range("D19").Value = "'10.75"
$ws.Range("E19").Value = "'  +0.11%  "

# Row 20
$ws.Range("D20").Value = "'324.19"
$ws.Range("E20").Value = "'  +1.38%  "

# Row 21
$ws.Range("E21").Value = "'  +0.56%  "

# Row 22
$ws.Range("D22").Value = "'6.08"
$ws.Range("E22").Value = "'  +6.21%  "

# Row 23
$ws.Range("E23").Value = "'  +0.00%  "

# Row 24
$ws.Range("D24").Value = "'63.74"
$ws.Range("E24").Value = "'  -0.40%  "

# Row 25
$ws.Range("D25").Value = "'0.407"
$ws.Range("E25").Value = "'  -0.75%  "

# Row 26
$ws.Range("E26").Value = "'  +1.08%  "

# Row 27
$ws.Range("E27").Value = "'  -0.80%  "

# Row 28
$ws.Range("D28").Value = "'7.38"
$ws.Range("E28").Value = "'  +1.07%  "

# Row 29
$ws.Range("D29").Value = "'0.0₃0748"
$ws.Range("E29").Value = "'  -0.33%  "

# Row 30
$ws.Range("D30").Value = "'168.58"
$ws.Range("E30").Value = "'  +1.41%  "

# Row 31
$ws.Range("D31").Value = "'1.71"
$ws.Range("E31").Value = "'  +0.93%  "

# Row 32
$ws.Range("B32").Value = "'Fetch.AI"
$ws.Range("C32").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.19"
$ws.Range("E32").Value = "'  +3.47%  "

# Row 33
$ws.Range("B33").Value = "'Aptos"
$ws.Range("C33").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "'6.29"
$ws.Range("E33").Value = "'  +0.10%  "

# Row 34
$ws.Range("E34").Value = "'  -0.03%  "

# Row 35
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = "'  -0.15%  "

# Row 36
$ws.Range("D36").Value = "'18.06"
$ws.Range("E36").Value = "'  +0.29%  "

# Row 37
$ws.Range("E37").Value = "'  -3.02%  "

# Row 38
$ws.Range("E38").Value = "'  -0.44%  "

# Row 39
$ws.Range("D39").Value = "'36.85"
$ws.Range("E39").Value = "'  +0.72%  "

# Row 40
$ws.Range("E40").Value = "'  -0.89%  "

# Row 41
$ws.Range("D41").Value = "'0.780"
$ws.Range("E41").Value = "'  -1.03%  "

# Row 42
$ws.Range("D42").Value = "'279.75"
$ws.Range("E42").Value = "'  +3.10%  "

# Row 43
$ws.Range("B43").Value = "'Filecoin"
$ws.Range("C43").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.44"
$ws.Range("E43").Value = "'  -0.19%  "

# Row 44
$ws.Range("B44").Value = "'RenderToken"
$ws.Range("C44").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'5.02"
$ws.Range("E44").Value = "'  +0.32%  "

# Row 45
$ws.Range("D45").Value = "'0.599"
$ws.Range("E45").Value = "'  +1.67%  "

# Row 46
$ws.Range("B46").Value = "'Aave"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'122.74"
$ws.Range("E46").Value = "'  -2.43%  "

# Row 47
$ws.Range("B47").Value = "'Stellar"
$ws.Range("C47").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.0921"
$ws.Range("E47").Value = "'  +1.84%  "

# Row 48
$ws.Range("E48").Value = "'  +2.44%  "

# Row 49
$ws.Range("D49").Value = "'17.76"
$ws.Range("E49").Value = "'  -0.24%  "

# Row 50
$ws.Range("E50").Value = "'  +0.28%  "

# Row 51
$ws.Range("D51").Value = "'17.04"
$ws.Range("E51").Value = "'  +0.35%  "
